$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("snapshot")

# Update scraped_at timestamps in column K (rows 2-38) on the "snapshot" sheet
$ws1.Range("K2").Value2 = "2025-12-09T07:01:30.888025+00:00"
$ws1.Range("K3").Value2 = "2025-12-09T07:01:30.888058+00:00"
$ws1.Range("K4").Value2 = "2025-12-09T07:01:30.888077+00:00"
$ws1.Range("K5").Value2 = "2025-12-09T07:01:33.225376+00:00"
$ws1.Range("K6").Value2 = "2025-12-09T07:01:33.225403+00:00"
$ws1.Range("K7").Value2 = "2025-12-09T07:01:35.520665+00:00"
$ws1.Range("K8").Value2 = "2025-12-09T07:01:38.247142+00:00"
$ws1.Range("K9").Value2 = "2025-12-09T07:01:40.580724+00:00"
$ws1.Range("K10").Value2 = "2025-12-09T07:01:42.962441+00:00"
$ws1.Range("K11").Value2 = "2025-12-09T07:01:47.659092+00:00"
$ws1.Range("K12").Value2 = "2025-12-09T07:01:47.659121+00:00"
$ws1.Range("K13").Value2 = "2025-12-09T07:01:49.939599+00:00"
$ws1.Range("K14").Value2 = "2025-12-09T07:01:52.182933+00:00"
$ws1.Range("K15").Value2 = "2025-12-09T07:01:54.906945+00:00"
$ws1.Range("K16").Value2 = "2025-12-09T07:01:57.251007+00:00"
$ws1.Range("K17").Value2 = "2025-12-09T07:01:57.251038+00:00"
$ws1.Range("K18").Value2 = "2025-12-09T07:01:59.993172+00:00"
$ws1.Range("K19").Value2 = "2025-12-09T07:01:59.993201+00:00"
$ws1.Range("K20").Value2 = "2025-12-09T07:01:59.993218+00:00"
$ws1.Range("K21").Value2 = "2025-12-09T07:02:02.792245+00:00"
$ws1.Range("K22").Value2 = "2025-12-09T07:02:02.792273+00:00"
$ws1.Range("K23").Value2 = "2025-12-09T07:02:02.792290+00:00"
$ws1.Range("K24").Value2 = "2025-12-09T07:02:02.792307+00:00"
$ws1.Range("K25").Value2 = "2025-12-09T07:02:02.792327+00:00"
$ws1.Range("K26").Value2 = "2025-12-09T07:02:05.577727+00:00"
$ws1.Range("K27").Value2 = "2025-12-09T07:02:05.577757+00:00"
$ws1.Range("K28").Value2 = "2025-12-09T07:02:05.577776+00:00"
$ws1.Range("K29").Value2 = "2025-12-09T07:02:07.882485+00:00"
$ws1.Range("K30").Value2 = "2025-12-09T07:02:12.440619+00:00"
$ws1.Range("K31").Value2 = "2025-12-09T07:02:12.440649+00:00"
$ws1.Range("K32").Value2 = "2025-12-09T07:02:12.440665+00:00"
$ws1.Range("K33").Value2 = "2025-12-09T07:02:14.769163+00:00"
$ws1.Range("K34").Value2 = "2025-12-09T07:02:14.769194+00:00"
$ws1.Range("K35").Value2 = "2025-12-09T07:02:17.715921+00:00"
$ws1.Range("K36").Value2 = "2025-12-09T07:02:17.715952+00:00"
$ws1.Range("K37").Value2 = "2025-12-09T07:02:20.025537+00:00"
$ws1.Range("K38").Value2 = "2025-12-09T07:02:20.025564+00:00"

# Remove the resolved "new_injured" entry (row 2) now that it has been moved/processed
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Rows.Item(2).Delete()
